# Auto-generated: Add 2026-01-30 data to M column across affected sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("M2").Value = 416
$ws.Range("M3").Value = 457
$ws.Range("M4").Value = 127
$ws.Range("M6").Value = 346
$ws.Range("M7").Value = 1372

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("M6").Value = 9
$ws.Range("M7").Value = 41
$ws.Range("M8").Value = 90
$ws.Range("M11").Value = 16
$ws.Range("M14").Value = 10
$ws.Range("M15").Value = 17
$ws.Range("M18").Value = 15
$ws.Range("M19").Value = 48
$ws.Range("M20").Value = 45
$ws.Range("M29").Value = 69
$ws.Range("M33").Value = 50
$ws.Range("M37").Value = 64
$ws.Range("M41").Value = 6
$ws.Range("M42").Value = 45
$ws.Range("M43").Value = 11
$ws.Range("M44").Value = 8
$ws.Range("M52").Value = 19
$ws.Range("M53").Value = 13
$ws.Range("M54").Value = 25
$ws.Range("M55").Value = 13
$ws.Range("M57").Value = 4
$ws.Range("M63").Value = 3
$ws.Range("M64").Value = 14
$ws.Range("M65").Value = 28
$ws.Range("M67").Value = 38
$ws.Range("M69").Value = 5
$ws.Range("M71").Value = 8
$ws.Range("M76").Value = 16
$ws.Range("M80").Value = 7
$ws.Range("M84").Value = 5
$ws.Range("M89").Value = 20
$ws.Range("M91").Value = 19
$ws.Range("M94").Value = 18
$ws.Range("M99").Value = 38
$ws.Range("M101").Value = 1372

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("M2").Value = 2
$ws.Range("M7").Value = 10

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("M2").Value = 16
$ws.Range("M7").Value = 41

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("M2").Value = 7
$ws.Range("M7").Value = 16

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("M3").Value = 6
$ws.Range("M7").Value = 20

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("M3").Value = 31
$ws.Range("M6").Value = 13

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("M2").Value = 6
$ws.Range("M3").Value = 4
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("M3").Value = 3
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("M2").Value = 5
$ws.Range("M7").Value = 13

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("M3").Value = 30
$ws.Range("M6").Value = 28
$ws.Range("M7").Value = 90

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("M2").Value = 11
$ws.Range("M3").Value = 17
$ws.Range("M7").Value = 50

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("M2").Value = 20
$ws.Range("M3").Value = 24
$ws.Range("M6").Value = 15
$ws.Range("M7").Value = 64

$ws = $wb.Worksheets.Item("New City")
$ws.Range("M2").Value = 9
$ws.Range("M7").Value = 28

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("M3").Value = 15
$ws.Range("M6").Value = 8
$ws.Range("M7").Value = 38

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("M3").Value = 12
$ws.Range("M7").Value = 38

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("M3").Value = 2
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("M2").Value = 4
$ws.Range("M7").Value = 25

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("M3").Value = 22
$ws.Range("M4").Value = 7
$ws.Range("M6").Value = 17
$ws.Range("M7").Value = 69

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("M4").Value = 2
$ws.Range("M7").Value = 48

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("M3").Value = 4
$ws.Range("M7").Value = 8

$ws = $wb.Worksheets.Item("River North")
$ws.Range("M3").Value = 3
$ws.Range("M4").Value = 2
$ws.Range("M7").Value = 16

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("M4").Value = 1
$ws.Range("M7").Value = 9

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("M2").Value = 3
$ws.Range("M7").Value = 6

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("M2").Value = 11
$ws.Range("M7").Value = 45

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("M2").Value = 6
$ws.Range("M7").Value = 13

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("M2").Value = 4
$ws.Range("M3").Value = 11
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("M2").Value = 4
$ws.Range("M7").Value = 14

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("M2").Value = 19
$ws.Range("M3").Value = 10
$ws.Range("M7").Value = 45

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("M6").Value = 5
$ws.Range("M7").Value = 15

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("M2").Value = 6
$ws.Range("M3").Value = 4
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 18

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("M3").Value = 5
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("M2").Value = 3
$ws.Range("M7").Value = 11

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("M3").Value = 3
$ws.Range("M7").Value = 4

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("M3").Value = 3
$ws.Range("M7").Value = 8

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("M2").Value = 1
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 7

Write-Host "Applied 2026-01-30 updates to $($wb.Worksheets.Count) total worksheets across 36 affected sheets."